$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.493.57'
$ws.Range('E2').Value = '  +4.58%  '
$ws.Range('D3').Value = '1.586.95'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -1.02%  '
$ws.Range('D5').Value = '213.86'
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('E6').Value = '  +0.74%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -1.11%  '
$ws.Range('D8').Value = '23.96'
$ws.Range('E8').Value = '  +8.55%  '
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').Value = '0.0887'
$ws.Range('E11').Value = '  +2.13%  '
$ws.Range('D12').Value = '1.811.25'
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').Value = '1.583.88'
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('D15').Value = '0.531'
$ws.Range('E15').Value = '  +1.77%  '
$ws.Range('D16').Value = '28.432.73'
$ws.Range('E16').Value = '  +4.48%  '
$ws.Range('D17').Value = '63.98'
$ws.Range('E17').Value = '  +2.65%  '
$ws.Range('D18').Value = '233.78'
$ws.Range('E18').Value = '  +8.08%  '
$ws.Range('D19').Value = '0.0₃0709'
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('D20').Value = '7.50'
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').Value = '0.997'
$ws.Range('E21').Value = '  -0.99%  '
$ws.Range('E22').Value = '  -0.58%  '
$ws.Range('D23').Value = '9.39'
$ws.Range('E23').Value = '  +1.43%  '
$ws.Range('D24').Value = '1.95'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').Value = '151.45'
$ws.Range('E25').Value = '  -1.60%  '
$ws.Range('D26').Value = '15.32'
$ws.Range('E26').Value = '  +1.47%  '
$ws.Range('D27').Value = '6.62'
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('E28').Value = '  +0.22%  '
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.96%  '
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('D31').Value = '0.0474'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').Value = '3.24'
$ws.Range('E32').Value = '  -0.46%  '
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('D34').Value = '1.417.51'
$ws.Range('E34').Value = '  -2.44%  '
$ws.Range('E35').Value = '  -1.36%  '
$ws.Range('E36').Value = '  -5.74%  '
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').Value = '2.54'
$ws.Range('E39').Value = '  +7.99%  '
$ws.Range('D40').Value = '0.543'
$ws.Range('E40').Value = '  +1.18%  '
$ws.Range('D41').Value = '0.813'
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  -0.98%  '
$ws.Range('D43').Value = '5.68'
$ws.Range('E43').Value = '  -2.34%  '
$ws.Range('D44').Value = '1.83'
$ws.Range('E44').Value = '  +5.47%  '
$ws.Range('E45').Value = '  -2.92%  '
$ws.Range('D46').Value = '64.41'
$ws.Range('E46').Value = '  -0.43%  '
$ws.Range('D47').Value = '1.722.66'
$ws.Range('E47').Value = '  +0.66%  '
$ws.Range('D48').Value = '87.48'
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('D49').Value = '0.0₆0107'
$ws.Range('E49').Value = '  +3.59%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').Value = '39.48'
$ws.Range('E51').Value = '  +16.38%  '
